$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "25.873.77"
$ws.Cells.Item(2, 5).Value = "  +0.30%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.638.27"
$ws.Cells.Item(3, 5).Value = "  +0.87%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.12%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "215.50"
$ws.Cells.Item(5, 5).Value = "  +0.19%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.41%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.05%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.84%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +1.13%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "20.22"
$ws.Cells.Item(10, 5).Value = "  +4.60%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0781"
$ws.Cells.Item(11, 5).Value = "  +0.44%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "Polkadot"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "4.27"
$ws.Cells.Item(12, 5).Value = "  +0.31%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.644.72"
$ws.Cells.Item(13, 5).Value = "  +1.25%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.863.85"
$ws.Cells.Item(14, 5).Value = "  +0.84%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.565"
$ws.Cells.Item(15, 5).Value = "  +1.67%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +1.95%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "63.31"
$ws.Cells.Item(17, 5).Value = "  -0.49%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "25.886.00"
$ws.Cells.Item(18, 5).Value = "  +0.32%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.07%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "194.04"
$ws.Cells.Item(20, 5).Value = "  +0.00%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "4.38"
$ws.Cells.Item(21, 5).Value = "  +1.16%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +1.79%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +4.11%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +0.12%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.76"
$ws.Cells.Item(25, 5).Value = "  -2.86%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "138.39"
$ws.Cells.Item(26, 5).Value = "  -2.09%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -3.76%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +1.65%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.56"
$ws.Cells.Item(29, 5).Value = "  +0.92%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.54%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +1.86%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.50%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.93%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +1.27%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +1.03%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.908"
$ws.Cells.Item(36, 5).Value = "  +1.53%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +1.93%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.551"
$ws.Cells.Item(38, 5).Value = "  +0.24%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "1.126.70"
$ws.Cells.Item(39, 5).Value = "  +0.03%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +0.64%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -0.07%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.49"
$ws.Cells.Item(42, 5).Value = "  -1.70%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "99.42"
$ws.Cells.Item(43, 5).Value = "  +2.47%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.801"
$ws.Cells.Item(44, 5).Value = "  +0.93%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "0.0₆0111"
$ws.Cells.Item(45, 5).Value = "  -0.85%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "55.52"
$ws.Cells.Item(46, 5).Value = "  +1.38%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -4.05%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.36%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.68"
$ws.Cells.Item(49, 5).Value = "  +1.24%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.00"
$ws.Cells.Item(50, 5).Value = "  -0.51%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +0.10%  "
